# Auto-generated edit script: updates crypto price/volume table (cols B-E, rows 2-51)
# matching the scraped coinranking.com data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.217.81"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.866.85"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = $ws.Range("D25").Style
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'0.7115"
$ws.Range("D5").Style = $ws.Range("D25").Style
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'241.70"
$ws.Range("D6").Style = $ws.Range("D25").Style
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = $ws.Range("D25").Style
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.3118"
$ws.Range("D8").Style = $ws.Range("D25").Style
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.07672"
$ws.Range("D9").Style = $ws.Range("D25").Style
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("D10").Value = "'24.74"
$ws.Range("D10").Style = $ws.Range("D25").Style
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").Value = "'0.08366"
$ws.Range("D11").Style = $ws.Range("D25").Style
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "1.865.53"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "'5.229"
$ws.Range("D13").Style = $ws.Range("D25").Style
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "'0.7119"
$ws.Range("D14").Style = $ws.Range("D25").Style
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "'91.38"
$ws.Range("D15").Style = $ws.Range("D25").Style
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "29.224.08"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "'5.942"
$ws.Range("D17").Style = $ws.Range("D25").Style
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "'243.71"
$ws.Range("D18").Style = $ws.Range("D25").Style
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'0.000007815"
$ws.Range("D19").Style = $ws.Range("D25").Style
$ws.Range("D20").Value = "2.114.45"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").Value = "'13.11"
$ws.Range("D21").Style = $ws.Range("D25").Style
$ws.Range("D22").Value = "'0.9987"
$ws.Range("D22").Style = $ws.Range("D25").Style
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'7.869"
$ws.Range("D23").Style = $ws.Range("D25").Style
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").Value = "'0.9996"
$ws.Range("D24").Style = $ws.Range("D25").Style
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").Value = "'163.28"
$ws.Range("D26").Style = $ws.Range("D25").Style
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'8.949"
$ws.Range("D27").Style = $ws.Range("D25").Style
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "'18.49"
$ws.Range("D28").Style = $ws.Range("D25").Style
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.501"
$ws.Range("D29").Style = $ws.Range("D25").Style
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.324"
$ws.Range("D30").Style = $ws.Range("D25").Style
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "'4.258"
$ws.Range("D32").Style = $ws.Range("D25").Style
$ws.Range("E32").Value = "  +3.33%  "
$ws.Range("D33").Value = "'0.05160"
$ws.Range("D33").Style = $ws.Range("D25").Style
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").Value = "'0.7979"
$ws.Range("D34").Style = $ws.Range("D25").Style
$ws.Range("E34").Value = "  +9.42%  "
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").Value = "'1.167"
$ws.Range("D36").Style = $ws.Range("D25").Style
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("D37").Value = "'2.684"
$ws.Range("D37").Style = $ws.Range("D25").Style
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'0.01852"
$ws.Range("D38").Style = $ws.Range("D25").Style
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("D39").Value = "'2.711"
$ws.Range("D39").Style = $ws.Range("D25").Style
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "1.164.41"
$ws.Range("E40").Value = "  -5.72%  "
$ws.Range("D41").Value = "'6.292"
$ws.Range("D41").Style = $ws.Range("D25").Style
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").Value = "'0.8967"
$ws.Range("D42").Style = $ws.Range("D25").Style
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("D43").Value = "'73.27"
$ws.Range("D43").Style = $ws.Range("D25").Style
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("D44").Value = "'0.9992"
$ws.Range("D44").Style = $ws.Range("D25").Style
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'103.08"
$ws.Range("D45").Style = $ws.Range("D25").Style
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "2.011.26"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").Value = "'1.781"
$ws.Range("D48").Style = $ws.Range("D25").Style
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").Value = "'9.346"
$ws.Range("D49").Style = $ws.Range("D25").Style
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'0.4295"
$ws.Range("D51").Style = $ws.Range("D25").Style
$ws.Range("E51").Value = "  -1.11%  "
